$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty score cells with their graded values.
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("J3").Value = 0

$ws.Range("B4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("L4").Value = 0

$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("L5").Value = 0

$ws.Range("B6").Value = 3
$ws.Range("D6").Value = 3
$ws.Range("H6").Value = 2
$ws.Range("J6").Value = 1

$ws.Range("B7").Value = 3
$ws.Range("D7").Value = 3
$ws.Range("F7").Value = 3
$ws.Range("H7").Value = 2
$ws.Range("J7").Value = 3

$ws.Range("B8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("H8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 1

$ws.Range("D10").Value = 0
$ws.Range("F10").Value = 3
$ws.Range("H10").Value = 1
$ws.Range("J10").Value = 1

$ws.Range("B11").Value = 3
$ws.Range("D11").Value = 3
$ws.Range("H11").Value = 2
$ws.Range("J11").Value = 2

$ws.Range("B12").Value = 2
$ws.Range("F12").Value = 3
$ws.Range("H12").Value = 1
$ws.Range("J12").Value = 0

$ws.Range("B13").Value = 3
$ws.Range("H13").Value = 3
$ws.Range("J13").Value = 2

$ws.Range("F14").Value = 1

$ws.Range("F15").Value = 1

$ws.Range("B16").Value = 2
$ws.Range("H16").Value = 2
$ws.Range("J16").Value = 3

$ws.Range("D18").Value = 2
$ws.Range("J18").Value = 2

$ws.Range("B19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("H19").Value = 0

$ws.Range("B20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 2
$ws.Range("H20").Value = 1
$ws.Range("J20").Value = 2

$ws.Range("B21").Value = 1
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0

$ws.Range("B22").Value = 2
$ws.Range("D22").Value = 3
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 1

# Freeze the header row and scroll so column K is leftmost visible,
# with the active selection at K15 (matches the saved view state).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("K15").Select()
$excel.ActiveWindow.ScrollColumn = 11
